$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.239.47'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '3.510.00'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'598.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = "'173.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").Value = "'7.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("D11").Value = "'0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '4.106.86'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").Value = "'29.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.55%  '
$ws.Range("D15").Value = '67.151.10'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '3.513.96'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = "'6.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = "'14.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("D20").Value = "'393.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = "'7.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = "'73.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = "'0.537"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = "'5.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").Value = "'10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = "'0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").Value = "'6.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.28%  '
$ws.Range("D31").Value = "'1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.09%  '
$ws.Range("D32").Value = "'2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").Value = "'23.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("D34").Value = "'7.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").Value = "'1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").Value = "'163.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").Value = "'0.877"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'6.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").Value = "'4.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").Value = "'27.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").Value = '2.815.96'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = "'0.0731"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").Value = "'26.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = "'2.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").Value = "'42.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("E47").Value = '  -2.29%  '
$ws.Range("D48").Value = "'343.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").Value = "'1.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").Value = "'33.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").Value = "'6.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
